$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "sluiten toegevoegd als bevestiging van projectlieders dat iets gesloten kan
# worden" -> insert a new "Sluiten" row right after "Bespreekpunten" (row 21),
# pushing the existing rows 22-35 down to 23-36.
$ws.Rows.Item(22).Insert()

# New row 22: mirrors the "Output"/"Hybrid" style rows around it (e.g. row 16 -
# "Verwacht resultaat") - appears in columns A, C and G but not D/E/F.
$ws.Range("A22").Value = "Sluiten"
$ws.Range("B22").Value = "Hybrid"
$ws.Range("C22").Value = "Sluiten"
$ws.Range("G22").Value = "Sluiten"

# Pick up the cell formatting (borders/fills) from the row below it (now row 23,
# the old row 22) so the new row matches the sheet's existing visual pattern.
$ws.Range("A23:G23").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# PasteSpecial(xlPasteFormats) does not touch values, so (re)write them now.
$ws.Range("A22").Value = "Sluiten"
$ws.Range("B22").Value = "Hybrid"
$ws.Range("C22").Value = "Sluiten"
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = ""
$ws.Range("G22").Value = "Sluiten"

# A handful of rows further down (now rows 28-36) carried slightly different,
# but visually identical, cell styles than the rest of column A/C. Normalize
# them to match the common style used across the sheet.
$ws.Range("A24").Copy()
$ws.Range("A31:A36").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C24").Copy()
$ws.Range("C28:C36").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A23").Select()
